$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 6: ALIGNMENT (D6) Left -> Right
$ws.Range("D6").Value = "Right"

# Row 7: GB (E7) "2" -> "9" (keep it text, like the existing "2" was, instead
# of letting it be auto-coerced to a number)
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "9"
$ws.Range("E7").Style = $ws.Range("D7").Style

# Rows 15-17 held the pOil / tOil / tWat texts that are now represented by
# icons instead, so remove those rows entirely (row 18, pFuel, shifts up to
# become the new row 15).
$ws.Range("B15:F17").EntireRow.Delete()
